$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStr = '"Medical complications included respiratory insufficiency in 12 cases, cardiac insufficiency in four, respiratory infection in three, acute renal failure in two, cardiac arrhythmia in one, acute gout in one, reactive depression in one, diabetic decompensation in one, acute pulmonary oedema in one and ventriculitis in one.", "Surgical complications included wound infection in nine cases (debridement and local treatment), paralytic ileus in six, intra-abdominal abscess in two (percutaneous drainage), partial dehiscence of the colostomy in two, perineal abscess in one (surgical drainage), cellulitis in one, evisceration in one (repaired with mesh) and haemoperitoneum in one (laparotomy and haemostasis)."'

# Populate the new row 28 with the extracted study data
$ws.Range("A28").Value = 588
$ws.Range("C28").Value = 33
$ws.Range("D28").Value = 33
$ws.Range("E28").Value = 19
$ws.Range("F28").Value = 96
$ws.Range("G28").Value = 64
$ws.Range("H28").Value = "UK"
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 1
$ws.Range("K28").Value = "UK"
$ws.Range("L28").Value = 0.455
$ws.Range("M28").Value = "UK"
$ws.Range("N28").Value = 0.545
$ws.Range("O28").Value = "UK"
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = "UK"
$ws.Range("R28").Value = 0.121
$ws.Range("S28").Value = "UK"
$ws.Range("T28").Value = "UK"
$ws.Range("U28").Value = "UK"
$ws.Range("V28").Value = "UK"
$ws.Range("W28").Value = "UK"
$ws.Range("X28").Value = "UK"
$ws.Range("Y28").Value = "UK"
$ws.Range("Z28").Value = 0.121
$ws.Range("AA28").Value = "UK"
$ws.Range("AB28").Value = 0.272
$ws.Range("AC28").Value = "UK"
$ws.Range("AD28").Value = "UK"
$ws.Range("AE28").Value = "UK"
$ws.Range("AF28").Value = "UK"
$ws.Range("AG28").Value = "UK"
$ws.Range("AH28").Value = "UK"
$ws.Range("AI28").Value = "UK"
$ws.Range("AJ28").Value = "UK"
$ws.Range("AK28").Value = "UK"
$ws.Range("AL28").Value = "UK"
$ws.Range("AM28").Value = "UK"
$ws.Range("AN28").Value = "UK"
$ws.Range("AO28").Value = 0
$ws.Range("AP28").Value = 1
$ws.Range("AQ28").Value = "UK"
$ws.Range("AR28").Value = "UK"
$ws.Range("AS28").Value = "UK"
$ws.Range("AT28").Value = "UK"
$ws.Range("AU28").Value = "UK"
$ws.Range("AV28").Value = "UK"
$ws.Range("AW28").Value = 33
$ws.Range("AX28").Value = 0
$ws.Range("AY28").Value = 0
$ws.Range("AZ28").Value = 0
$ws.Range("BA28").Value = 0
$ws.Range("BB28").Value = "UK"
$ws.Range("BC28").Value = 0.091
$ws.Range("BD28").Value = "UK"
$ws.Range("BE28").Value = 0.515
$ws.Range("BF28").Value = "UK"
$ws.Range("BG28").Value = 0.909
$ws.Range("BH28").Value = "UK"
$ws.Range("BI28").Value = "UK"
$ws.Range("BJ28").Value = "UK"
$ws.Range("BK28").Value = 0.121
$ws.Range("BL28").Value = "UK"
$ws.Range("BM28").Value = 0.879
$ws.Range("BN28").Value = "UK"
$ws.Range("BO28").Value = 0.8181818182
$ws.Range("BP28").Value = 0
$ws.Range("BQ28").Value = 0
$ws.Range("BR28").Value = "UK"
$ws.Range("BS28").Value = 1
$ws.Range("BT28").Value = "UK"
$ws.Range("BU28").Value = 0.061
$ws.Range("BV28").Value = "UK"
$ws.Range("BW28").Value = 1
$ws.Range("BX28").Value = "UK"
$ws.Range("BY28").Value = 1
$ws.Range("BZ28").Value = "UK"
$ws.Range("CA28").Value = "UK"
$ws.Range("CB28").Value = "UK"
$ws.Range("CC28").Value = 0.4545454545
$ws.Range("CD28").Value = $newStr
$ws.Range("CE28").Value = 0
$ws.Range("CF28").Value = 0
$ws.Range("CG28").Value = 0
$ws.Range("CH28").Value = 0

# Match the author's final on-screen selection after the edit.
$ws.Range("F28").Select()
